$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record as row 26, pushing the existing rows 26-37 (and their
# data) down to 27-38.
$ws.Rows("26:26").Insert()

# Fill in the new row 26 with the new weekly price record.
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 44460
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112031
$ws.Cells.Item(26, 7).Value = "Poroto verde"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 1200
$ws.Cells.Item(26, 11).Value = 1400
$ws.Cells.Item(26, 12).Value = 1500
$ws.Cells.Item(26, 13).Value = 1450
$ws.Cells.Item(26, 14).Value = "`$/kilo"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 1450
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
